$wb = $excel.ActiveWorkbook

# "montage financier" sheet: selection moves from D25 to C22:D33,
# and it loses tab-selected status to "spécification".
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C22:D33").Select()

# "spécification" sheet: new domotique / droit-de-passage notes.
$ws3 = $wb.Worksheets.Item(3)

# Write cells in this exact order so new shared strings land at the
# indices the target workbook expects (41..44).
$ws3.Range("A18").Value = "http://www.digitaltrends.com/home/zigbee-vs-zwave-vs-insteon-home-automation-protocols-explained/#ixzz3xL6rNmK2"
$ws3.Range("A17").Value = "protocoles"
$ws3.Range("A16").Value = "Domotique"
$ws3.Range("A21").Value = "Logiciel"
$ws3.Range("A22").Value = "openhab"

# Activate "spécification" last so it becomes the active tab/sheet,
# and select A22 as the final active cell on it.
$ws3.Activate()
$ws3.Range("A22").Select()
